$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.019745111465454
$ws.Range("B1").Value = 3.432794570922852
$ws.Range("C1").Value = 3.014929294586182
$ws.Range("D1").Value = 3.277546644210815
$ws.Range("E1").Value = 2.028753280639648
